$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("D11")
Write-Host ($cell.Value2)
$ws.Hyperlinks.Add($cell, $cell.Value2) | Out-Null
Write-Host "done"
